$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update designator lists / quantities for parts whose populated count changed ---

# 0.1uF 0402 caps: remove C49
$ws.Range("C2").Value = "C1, C3, C5, C7, C9, C11, C23, C25, C27, C44, C45, C46"
$ws.Range("B2").Value = 12

# 10uF 0603 caps: remove C38, C39
$ws.Range("C3").Value = "C2, C4, C6, C8, C10, C12, C15, C17, C18, C19, C20, C21, C24, C26, C28, C42, C43"
$ws.Range("B3").Value = 17

# 600R/100MHz ferrite beads: remove FB6, FB9
$ws.Range("C12").Value = "FB1, FB2, FB3, FB4, FB5, FB7, FB8"
$ws.Range("B12").Value = 7

# 45.3K 0402 resistors: remove R37
$ws.Range("C25").Value = "R8"
$ws.Range("B25").Value = 1

# --- Remove the 1.5V LDO regulator (U8, TLV70015DCK) row entirely ---
$ws.Rows.Item(36).Delete()

# --- Page setup tweak recorded alongside this revision ---
$ws.PageSetup.PaperSize = 9
